# Replace the "Cluster Name" / "Active Cases" outbreak table contents
# with the updated data set (renamed header, refreshed rows, new case counts).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "Cluster Name"
$ws.Cells.Item(1, 2).Value = "Activecases"

# Data rows 2-38: cluster name (col A) + active case count (col B)
$rows = @(
    @("202111 Holy Rosary Primary School White Hills", 25),
    @("3323 Villa Maria Catholic Homes St Bernadette's Aged Care Sunshine North", 15),
    @("3376 Royal Freemasons Coppin Centre Melbourne", 11),
    @("3601 Baptcare Westhaven community outbreak", 30),
    @("3653 Fronditha Thalpori St Albans Aged Care", 34),
    @("44121 Wallan Primary School Wallan", 20),
    @("44165 Greenvale Primary School", 37),
    @("44234 Lucknow Primary School Bairnsdale", 36),
    @("44395 Buln Buln Primary School", 11),
    @("44667 Beaumaris Primary School Beaumaris", 22),
    @("44811 Dandenong North Primary School Dandenong", 27),
    @("44820 Bourchier Street Primary School Shepparton", 10),
    @("44853 St Albans North Primary School", 12),
    @("45158 Rowellyn Park Primary School Carrum Downs", 16),
    @("45249 Creekside K-9 College Caroline Springs", 16),
    @("4574 Village Glen Aged Care Residences Mornington", 18),
    @("45812 St Mary's Primary School Hampton", 11),
    @("45836 St Joseph's Primary School Sorrento", 10),
    @("46037 Nazareth Catholic Primary School Grovedale", 27),
    @("46050 Our Lady's Catholic Primary School Craigieburn", 40),
    @("46125 Our Lady of the Southern Cross PrimarySchool Manor Lakes Outbreak", 29),
    @("46140 Caulfield Grammar School Malvern Campus Glen Iris", 10),
    @("46190 Haileybury Brighton East", 13),
    @("46276 Hillcrest Christian College Clyde North", 16),
    @("46328 Ilim College Dallas Main Campus DallasOct", 28),
    @("46328 Ilim College Dallas Primary Campus Inverloch Cres Tier 1A Dallas", 10),
    @("46376 Yesodei HaTorah College Elwood", 12),
    @("46390 Al Siraat College Epping", 42),
    @("50395 St Francis of Assisi Catholic Primary School Tarneit", 12),
    @("50681 Broadmeadows Special Developmental School Broadmeadows", 12),
    @("52380 Al Iman College Melton South", 14),
    @("52473 John Henry Primary School Pakenham", 18),
    @("Alfred Health Caulfield Hospital Caulfield", 13),
    @("Hamilton Country Music Festival Hamilton GolfClub Hamilton", 19),
    @("Melton Willows Melton", 11),
    @("Sirius College Shepparton Campus Shepparton", 11),
    @("The Sphinx Hotel North Geelong", 10),
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
